# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# This updates the K column values (G2:G16) to the recomputed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 3
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 0
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
